# Add visualisation config to simple chart survey for coverage
#
# The "Test Chart" worksheet holds the simpleChart survey's question rows.
# Row 4 is the "testchartcode2" (Number / Size) question; give it a
# validationCriteria (column K) and a visualisationConfig (column L) value
# so the importer test fixture exercises chart visualisation config parsing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Chart")

# --- Row 4: add validationCriteria (K4) and visualisationConfig (L4) ---
$k4 = $ws.Cells.Item(4, 11)   # K4
$k4.Value = '{"min": 0, "max": 300, "normalRange": {"min": 90, "max": 120}}'
$k4.Font.Name = "Calibri"
$k4.Font.Size = 12

# Copy K4's freshly created formatting onto L4, then set L4's own value so
# both new cells share a single new style instead of each spawning one.
$k4.Copy()
$l4 = $ws.Cells.Item(4, 12)   # L4
$l4.PasteSpecial(-4122)       # xlPasteFormats
$l4.Value = '{"yAxis":{"graphRange":{"min":40,"max":240},"interval":10}}'

# --- Row 2: re-normalise A2's style (drop the now-superfluous alignment
# override) to match the plain style used by its sibling cells B2:D2 ---
$b2 = $ws.Cells.Item(2, 2)
$b2.Copy()
$a2 = $ws.Cells.Item(2, 1)
$a2.PasteSpecial(-4122)       # xlPasteFormats
